# Disponibilidad.xlsx - "Actualizar 02-05-2021 09-18-14" automatic update
# 1) A tiny floating point re-stamp of the previous run's timestamp
#    (rows 464-477, column D).
# 2) Append a brand-new 14-row availability check block (rows 478-491)
#    with its own timestamp, replicating the existing per-row layout:
#    A=service name, B=URL (hyperlinked), C="Disponible", D=timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-stamp the previous block's cached timestamp (rows 464-477).
# ---------------------------------------------------------------------
for ($r = 464; $r -le 477; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.36656471065
}

# ---------------------------------------------------------------------
# 2) Append the new availability block (rows 478-491).
# ---------------------------------------------------------------------
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# The MapStore row (index 8, 0-based) carries a "#/" fragment in the
# displayed URL / a location="/" on the hyperlink relationship.
$subAddress = @($null,$null,$null,$null,$null,$null,$null,$null,"/",$null,$null,$null,$null,$null)

$timestamp = 44232.38761723032
$startRow = 478

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $startRow + $i
    $displayUrl = $urls[$i]
    if ($subAddress[$i]) {
        $displayUrl = $urls[$i] + "#" + $subAddress[$i]
    }

    $ws.Cells.Item($row, 1).Value = $names[$i]

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $displayUrl
    if ($subAddress[$i]) {
        $ws.Hyperlinks.Add($bCell, $urls[$i], $subAddress[$i])
    } else {
        $ws.Hyperlinks.Add($bCell, $urls[$i])
    }
    $bCell.Style = "Hyperlink"

    $ws.Cells.Item($row, 3).Value = "Disponible"

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $timestamp
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
